# Cover letter revision ("submission end of march 2021")
# Applies the textual edits described by the commit:
#  - Shortened/retitled manuscript title
#  - Expanded description of the study's scope ("model system" / broader implications)
#  - Trimmed claim about body size being the key generalizable trait
#  - Reworded sentence about predator species identity driving prey-size relationships
#  - Emphasized novelty of the interaction data
#  - Re-worked sentence introducing terrestrial invertebrates' ecological importance
#  - Specified "ecosystem-regulating" taxa
#  - Added a new closing paragraph on originality/conflicts of interest/prior publication

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Predator-prey interaction outcomes in terrestrial invertebrates are determined by predator body size and species identity, but not hunting traits as inferred from diet DNA metabarcoding data",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Predator-prey interactions of terrestrial invertebrates are determined by predator body size and species identity",
    2) | Out-Null

$d.Content.Find.Execute(
    "by providing an exploration of mechanistic and generalizable patterns of predator-prey interactions, linking interaction outcomes to predator traits. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "by providing an exploration of mechanistic and generalizable patterns of community predator-prey interactions, linking interaction outcomes to predator traits in a model system with broader implications across terrestrial invertebrate taxa globally. ",
    2) | Out-Null

$d.Content.Find.Execute(
    " suggest that generalizable predator traits, specifically body size, are key to shaping predator-prey interactions for terrestrial invertebrates, predicting the body size of prey ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " suggest that generalizable predator traits are key to shaping predator-prey interactions for terrestrial invertebrates, predicting the body size of prey ",
    2) | Out-Null

$d.Content.Find.Execute(
    ". Instead, predator species identity shaped the size of prey that predator individuals consumed. By using ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ". Instead, these relationships were driven by predator species identity. By using ",
    2) | Out-Null

$d.Content.Find.Execute(
    "DNA metabarcoding data, our study provides valuable empirical interaction data for a set of organisms (terrestrial invertebrates including spiders, insects, and centipedes) for which we have limited observed knowledge",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "DNA metabarcoding data, our study provides novel and highly valuable interaction data for a set of organisms (terrestrial invertebrates including spiders, insects, and centipedes) for which we have limited observed knowledge",
    2) | Out-Null

$d.Content.Find.Execute(
    "However, while we have limited knowledge of the interactions in this set of consumers, terrestrial invertebrates are the most abundant and diverse taxa on earth, so their functional roles in ecosystems are likely to strongly impact community structure and ecosystem dynamics. Therefore, w",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Terrestrial invertebrates are the most abundant and diverse taxa on earth, so their functional roles in ecosystems are likely to strongly impact community structure and ecosystem dynamics. Therefore, w",
    2) | Out-Null

$d.Content.Find.Execute(
    " among taxa for which we have had limited data",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " among ecosystem-regulating taxa for which we have had limited data",
    2) | Out-Null

$d.Content.Find.Execute(
    "We look forward to hearing your thoughts on this work. Please feel free to contact me with any questions about this material. On behalf of my co-authors, I thank you for your consideration of our submission.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The work in this manuscript is all original research carried out by the authors and all authors agree to its content. We have no conflicts of interest to report and have received appropriate approvals to conduct this research. The results in this manuscript have not been submitted for publication elsewhere, nor are they previously published. We look forward to hearing your thoughts on this work. Please feel free to contact me with any questions about this material. On behalf of my co-authors, I thank you for your consideration of our submission. ",
    2) | Out-Null
